# Update shadow model values according to D1/L1 = D2/L2 --> D2 = ...
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 8362.283772317136
$ws.Range("D5").Value = 8362.283772317136

$ws.Range("D9").Value = 12009.13023996808
$ws.Range("D10").Value = 12009.13023996808

$ws.Range("D14").Value = 11855.31622768284
$ws.Range("D15").Value = 11855.31622768284
